# Update gh-pages output data (江西-漫展信息)
# The workbook has 4 sheets: 展览(1), 演出(2), 本地生活(3), 全部类型(4)
# Sheets 1 (展览) and 4 (全部类型) carry the same event-listing table and
# both receive identical updates to the "想去人数" (F) and "最低票价" (G)
# columns, refreshed from the upstream scrape.

$wb = $excel.ActiveWorkbook

foreach ($sheetIndex in 1, 4) {
    $ws = $wb.Worksheets.Item($sheetIndex)

    $ws.Range("F2").Value = 2093

    # Ticket for this event is no longer sellable; replaced numeric price
    # with the Chinese text marker used elsewhere in the sheet.
    $ws.Range("G3").Value = "不可售"

    $ws.Range("F6").Value = 1751
    $ws.Range("F8").Value = 735
    $ws.Range("F15").Value = 146
    $ws.Range("F17").Value = 138
    $ws.Range("F18").Value = 4071
    $ws.Range("F21").Value = 452
    $ws.Range("F22").Value = 391
    $ws.Range("F23").Value = 954
    $ws.Range("F24").Value = 942
    $ws.Range("F26").Value = 25
    $ws.Range("F28").Value = 1843
    $ws.Range("F29").Value = 46
}
